# This script applies the updated NATMI TPM-derived statistics for the
# Tctn1 -> Tmem67 ligand-receptor pair sheet (all 4x4 sending/target
# cluster combinations). Values come from re-running the analysis with
# the new TPM expression matrix; written here as exact literals so the
# workbook matches the recomputed output byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G"="2.792060333333334"; "H"="8.376181000000001"; "I"="0.1656462876761766"; "J"="0.1656462876761766"; "M"="6.117613666666667"; "N"="18.352841"; "O"="0.3472444214548092"; "P"="0.3472444214548092"; "Q"="17.08074645335789"; "R"="153.726718080221"; "S"="0.05751974933025082"; "T"="0.05751974933025084" }
    3 = @{ "G"="2.792060333333334"; "H"="8.376181000000001"; "I"="0.1656462876761766"; "J"="0.1656462876761766"; "O"="0.2834964117229506"; "P"="0.2834964117229506"; "Q"="13.94501978977556"; "R"="125.50517810798"; "S"="0.04696012817142368"; "T"="0.04696012817142369" }
    4 = @{ "G"="2.792060333333334"; "H"="8.376181000000001"; "I"="0.1656462876761766"; "J"="0.1656462876761766"; "M"="3.555573333333333"; "N"="10.66672"; "O"="0.2018193812729289"; "P"="0.2018193812729289"; "Q"="9.927375266257778"; "R"="89.34637739632001"; "S"="0.03343063128896354"; "T"="0.03343063128896356" }
    5 = @{ "G"="2.792060333333334"; "H"="8.376181000000001"; "I"="0.1656462876761766"; "J"="0.1656462876761766"; "M"="2.949887333333333"; "N"="8.849661999999999"; "O"="0.1674397855493113"; "P"="0.1674397855493114"; "Q"="8.236263411202444"; "R"="74.12637070082199"; "S"="0.02773577888553854"; "T"="0.02773577888553855" }
    6 = @{ "I"="0.3370067148381872"; "J"="0.3370067148381872"; "M"="6.117613666666667"; "N"="18.352841"; "O"="0.3472444214548092"; "P"="0.3472444214548092"; "Q"="34.75071086702078"; "R"="312.756397803187"; "S"="0.1170237017203722"; "T"="0.1170237017203722" }
    7 = @{ "I"="0.3370067148381872"; "J"="0.3370067148381872"; "O"="0.2834964117229506"; "P"="0.2834964117229506"; "S"="0.09554019438316573"; "T"="0.09554019438316574" }
    8 = @{ "I"="0.3370067148381872"; "J"="0.3370067148381872"; "M"="3.555573333333333"; "N"="10.66672"; "O"="0.2018193812729289"; "P"="0.2018193812729289"; "Q"="20.19720557811556"; "R"="181.77485020304"; "S"="0.06801448667346534"; "T"="0.06801448667346534" }
    9 = @{ "I"="0.3370067148381872"; "J"="0.3370067148381872"; "M"="2.949887333333333"; "N"="8.849661999999999"; "O"="0.1674397855493113"; "P"="0.1674397855493114"; "Q"="16.75664522091489"; "R"="150.809806988234"; "S"="0.05642833206118399"; "T"="0.05642833206118399" }
    10 = @{ "G"="4.041018999999999"; "H"="12.123057"; "I"="0.2397440298074607"; "J"="0.2397440298074608"; "M"="6.117613666666667"; "N"="18.352841"; "O"="0.3472444214548092"; "P"="0.3472444214548092"; "Q"="24.72139306165966"; "R"="222.492537554937"; "S"="0.08324977692773622"; "T"="0.08324977692773625" }
    11 = @{ "G"="4.041018999999999"; "H"="12.123057"; "I"="0.2397440298074607"; "J"="0.2397440298074608"; "O"="0.2834964117229506"; "P"="0.2834964117229506"; "Q"="20.18297715600666"; "R"="181.64679440406"; "S"="0.06796657218241522"; "T"="0.06796657218241525" }
    12 = @{ "G"="4.041018999999999"; "H"="12.123057"; "I"="0.2397440298074607"; "J"="0.2397440298074608"; "M"="3.555573333333333"; "N"="10.66672"; "O"="0.2018193812729289"; "P"="0.2018193812729289"; "Q"="14.36813939589333"; "R"="129.31325456304"; "S"="0.04838499175962034"; "T"="0.04838499175962036" }
    13 = @{ "G"="4.041018999999999"; "H"="12.123057"; "I"="0.2397440298074607"; "J"="0.2397440298074608"; "M"="2.949887333333333"; "N"="8.849661999999999"; "O"="0.1674397855493113"; "P"="0.1674397855493114"; "Q"="11.92055076185933"; "R"="107.284956856734"; "S"="0.04014268893768892"; "T"="0.04014268893768894" }
    14 = @{ "G"="4.342041333333333"; "H"="13.026124"; "I"="0.2576029676781755"; "J"="0.2576029676781755"; "M"="6.117613666666667"; "N"="18.352841"; "O"="0.3472444214548092"; "P"="0.3472444214548092"; "Q"="26.56293140203156"; "R"="239.066382618284"; "S"="0.08945119347644996"; "T"="0.08945119347644997" }
    15 = @{ "G"="4.342041333333333"; "H"="13.026124"; "I"="0.2576029676781755"; "J"="0.2576029676781755"; "O"="0.2834964117229506"; "P"="0.2834964117229506"; "Q"="21.68644122710222"; "R"="195.17797104392"; "S"="0.07302951698594598"; "T"="0.07302951698594599" }
    16 = @{ "G"="4.342041333333333"; "H"="13.026124"; "I"="0.2576029676781755"; "J"="0.2576029676781755"; "M"="3.555573333333333"; "N"="10.66672"; "O"="0.2018193812729289"; "P"="0.2018193812729289"; "Q"="15.43844637703111"; "R"="138.94601739328"; "S"="0.05198927155087969"; "T"="0.05198927155087969" }
    17 = @{ "G"="4.342041333333333"; "H"="13.026124"; "I"="0.2576029676781755"; "J"="0.2576029676781755"; "M"="2.949887333333333"; "N"="8.849661999999999"; "O"="0.1674397855493113"; "P"="0.1674397855493114"; "Q"="12.80853273000978"; "R"="115.276794570088"; "S"="0.04313298566489989"; "T"="0.0431329856648999" }
}

foreach ($r in $updates.Keys) {
    $cols = $updates[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = [double]$cols[$col]
    }
}
